# Fixed formatting: add borders and clean up header row, and drop the
# unused "No. of Sites ..." / accomplishment tracking columns (X..AG),
# keeping the last "Status as of ..." column which slides left into X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused columns X through AG (10 columns). Excel shifts
# the remaining "Status as of July 4, 2025" column (AH) left into X.
$ws.Range("X1:AG2").EntireColumn.Delete()

# The cell that lands in X1 (old AH1) carries no special formatting.
# Give the header cell the same bold font + border used by the rest of
# row 1, but without the centered/top alignment applied to the others.
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").HorizontalAlignment = 1
$ws.Range("X1").VerticalAlignment = -4107

# Add a border around every data cell in row 2 (A2:X2), matching the
# border already used on the header row.
$ws.Range("A2:X2").Borders.LineStyle = 1
$ws.Range("A2:X2").Borders.Weight = 2
